# Add a new tracking row (row 14) to the Model training tracking sheet,
# reproducing a training run with a lower min cluster size (F=5) and the
# "3000s ca " training time, analogous to the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "regular"
$ws.Range("B14").Value = "full random"
$ws.Range("C14").Value = "sentences"
$ws.Range("D14").Value = 3000
$ws.Range("E14").Value = 200
$ws.Range("F14").Value = 5
$ws.Range("H14").Value = "3000s ca "
$ws.Range("I14").Value = "auto"
$ws.Range("J14").Value = "yes"
$ws.Range("K14").Value = 115
$ws.Range("L14").Value = "bad"

# Match the saved selection state observed in the edited workbook
$ws.Range("M15").Select()
